$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 1
    7 = 1
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 3
    18 = 0
    19 = 2
    20 = 2
    21 = 3
    22 = 0
    23 = 0
    24 = 0
    25 = 3
    26 = 0
    27 = 0
    29 = 0
    30 = 1
    31 = 2
    32 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 1
    39 = 3
    40 = 1
    41 = 0
    42 = 0
    43 = 3
    44 = 0
    45 = 4
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 2
    52 = 3
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 0
    58 = 1
    59 = 2
    60 = 3
    61 = 1
    62 = 1
    63 = 1
    64 = 2
    65 = 3
    66 = 2
    67 = 1
    68 = 1
    69 = 3
    70 = 1
    71 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

